$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows for the second classification method (Least Distance) applied
# to the "Times New Roman" set -- proofreading results from multiple
# classification algorithms on the same image set.
$ws.Range("A6").Value = "Times New Roman"
$ws.Range("B6").Value = "82/92"
$ws.Range("C6").Value = 0.891
$ws.Range("D6").Value = "82/92 (89.1%)"
$ws.Range("E6").Value = "Least Distance"

$ws.Range("A7").Value = "Times New Roman"
$ws.Range("B7").Value = "88/92"
$ws.Range("C7").Value = 0.957
$ws.Range("D7").Value = "90/92 (97.8%)"
$ws.Range("E7").Value = "Least Distance"

# Reuse the existing percentage style from the row above instead of
# minting a brand-new number format for the recognition-rate column.
$ws.Range("C5").Copy()
$ws.Range("C6:C7").PasteSpecial(-4122)

# Summary row with the average recognition rate across all methods.
$ws.Range("A17").Value = "Average Recognition"
$ws.Range("C17").Formula = "=AVERAGE(C2:C7)"
# AVERAGE() would otherwise inherit the percentage format from the cells
# it references -- put it back to the default/general style.
$ws.Range("C17").Style = "Normal"

$ws.Range("A17").Select()
